$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: update existing row 37 ---
# Date moves from 2022-03-23 (44643) to 2022-05-13 (44694)
$ws.Range("D37").Value = 44694
# Volume changes from 100 to 80
$ws.Range("M37").Value = 80

# --- Step 2: insert two new rows before the old row 38 ---
# This pushes the old row 38 (date 44644) down to row 40, keeping its values intact.
$ws.Rows("38:39").Insert()

# --- Step 3: populate new row 38 ("Segunda" quality, same new date) ---
$ws.Range("A38").Value = 6
$ws.Range("B38").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C38").Value = "Metropolitana"
$ws.Range("D38").Value = 44694
$ws.Range("E38").Value = 13
$ws.Range("F38").Value = "Fruta"
$ws.Range("G38").Value = 100101
$ws.Range("H38").Value = "Berries"
$ws.Range("I38").Value = 100101006
$ws.Range("J38").Value = "Higo"
$ws.Range("K38").Value = "Sin especificar"
$ws.Range("L38").Value = "Segunda"
$ws.Range("M38").Value = 75
$ws.Range("N38").Value = 12000
$ws.Range("O38").Value = 12000
$ws.Range("P38").Value = 12000
$ws.Range("Q38").Value = "$/bandeja 7 kilos"
$ws.Range("R38").Value = "Región Metropolitana"
$ws.Range("S38").Value = 1714
$ws.Range("T38").Value = 7

# --- Step 4: populate new row 39 (restore the original row-37 record) ---
$ws.Range("A39").Value = 6
$ws.Range("B39").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C39").Value = "Metropolitana"
$ws.Range("D39").Value = 44643
$ws.Range("E39").Value = 13
$ws.Range("F39").Value = "Fruta"
$ws.Range("G39").Value = 100101
$ws.Range("H39").Value = "Berries"
$ws.Range("I39").Value = 100101006
$ws.Range("J39").Value = "Higo"
$ws.Range("K39").Value = "Sin especificar"
$ws.Range("L39").Value = "Primera"
$ws.Range("M39").Value = 100
$ws.Range("N39").Value = 15000
$ws.Range("O39").Value = 15000
$ws.Range("P39").Value = 15000
$ws.Range("Q39").Value = "$/bandeja 7 kilos"
$ws.Range("R39").Value = "Región Metropolitana"
$ws.Range("S39").Value = 2143
$ws.Range("T39").Value = 7

# --- Step 5: row 40 already holds the old row-38 values (date 44644, etc.) ---
# verified by the prior Insert step; style for D38/D39 was inherited as the date format.

# Make sure the date cells keep the date number format used by the rest of column D.
$ws.Range("D38").NumberFormat = $ws.Range("D40").NumberFormat
$ws.Range("D39").NumberFormat = $ws.Range("D40").NumberFormat
